$wb = $excel.ActiveWorkbook

# --- Sheets ---
$active = $wb.Worksheets.Item("Active")
$inactive = $wb.Worksheets.Item("Inactive")
$config = $wb.Worksheets.Item("Config")

# --- 1. The old "Active" task (Id 64, "edit palette: edit color in place") is
#        now finished, so it gets replaced in-place on the Active sheet by a
#        brand-new task that takes the next Max Id ---
$oldId = $active.Range("A5").Value2
$oldTitle = $active.Range("B5").Value2
$oldCategory = $active.Range("D5").Value2
$oldCreated = $active.Range("E5").Value2

$newMaxId = $config.Range("F2").Value2 + 1

$active.Range("A5").Value2 = $newMaxId
$active.Range("B5").Value2 = "edit palette: add new color, but start from existing color"

# --- 2. Move the completed task down into "Inactive" as a new row 2,
#        pushing the existing rows down one ---
$inactive.Rows.Item(2).Insert()
$inactive.Rows.Item(2).ClearFormats()

$inactive.Range("A2").Value2 = $oldId
$inactive.Range("B2").Value2 = $oldTitle
$inactive.Range("C2").Value2 = "Done"
$inactive.Range("D2").Value2 = $oldCategory

# the "Created"/"Done" columns hold plain text dates (e.g. "8/24/2018"), not
# real date serials -- force text formatting so Excel doesn't auto-convert
$inactive.Range("E2:F2").NumberFormat = "@"
$inactive.Range("E2").Value2 = $oldCreated
$inactive.Range("F2").Value2 = $oldCreated

# --- 3. Bump "Max Id" on the Config sheet ---
$config.Range("F2").Value2 = $newMaxId
